$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark (bookmarkStart/bookmarkEnd) from the date paragraph.
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Replace the "所属" paragraph: merge the two fitText runs ("所" + "属") into a single
#    plain run "所属" and change the paragraph indent from left/firstLine to firstLineChars/firstLine.
$pShozoku = $d.Paragraphs.Item(22)
$xmlShozoku = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14 w15"><w:body><w:p w:rsidR="00704C43" w:rsidRDefault="003E7264"><w:pPr><w:ind w:firstLineChars="1300" w:firstLine="4680"/><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>所属</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pShozoku.Range.InsertXML($xmlShozoku)

# 3. Replace the big tab-filled paragraph (now at index 23) that contains the
#    "　　　" + tabs + fitText("氏"+"名") + tabs run sequence with four new paragraphs:
#      a) an empty paragraph
#      b) a paragraph containing "氏名" + "（自署）"
#      c) a paragraph with five tab runs
#      d) a paragraph with four tab runs (same ind as before)
$pBig = $d.Paragraphs.Item(23)
$xmlBig = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14 w15"><w:body><w:p w:rsidR="003E7264" w:rsidRDefault="003E7264"><w:pPr><w:tabs><w:tab w:val="left" w:pos="4140"/></w:tabs><w:spacing w:line="480" w:lineRule="exact"/><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="4140"/></w:tabs><w:spacing w:line="480" w:lineRule="exact"/><w:ind w:firstLineChars="1300" w:firstLine="4680"/><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>氏名</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="36"/></w:rPr><w:t>（自署）</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="4140"/></w:tabs><w:spacing w:line="480" w:lineRule="exact"/><w:ind w:firstLineChars="1200" w:firstLine="4320"/><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:tab/></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="4140"/></w:tabs><w:spacing w:line="480" w:lineRule="exact"/><w:ind w:left="4140"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="ＭＳ 明朝" w:eastAsia="ＭＳ 明朝" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:tab/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pBig.Range.InsertXML($xmlBig)
